# Applies the historico_dividendos.xlsx update:
#  - Divide a batch of mis-scaled dividend values (column C, rows 57-81) by 10
#  - Add an AutoFilter over the data range A1:C193 (+ the _FilterDatabase defined name)
#  - Move the sheet selection to B115 (scrolled back to the top)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("df_dividendos")

# --- Fix column C values that were off by a factor of 10 ---------------
# (use the literal decimal targets rather than "/10" so the stored double
#  bit-pattern matches Excel's own literal parse, not a divided-at-runtime
#  value which can differ in the last bit, e.g. 0.82/10 != 0.082)
$fixedValues = @{
    57 = 0.09
    58 = 0.09
    59 = 0.09
    62 = 0.09
    63 = 0.09
    64 = 0.09
    66 = 0.09
    67 = 0.093
    68 = 0.093
    69 = 0.094
    70 = 0.091
    71 = 0.091
    72 = 0.089
    73 = 0.089
    74 = 0.089
    75 = 0.091
    76 = 0.089
    77 = 0.089
    78 = 0.082
    79 = 0.09
    80 = 0.088
    81 = 0.092
}

foreach ($r in $fixedValues.Keys) {
    $ws.Cells.Item($r, 3).Value = $fixedValues[$r]
}

# --- Add AutoFilter across the used range + its defined name -----------
$ws.Range("A1:C193").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=df_dividendos!`$A`$1:`$C`$193", $false)
$fdb.Visible = $false

# --- Reset the view: selection on B115, scrolled to top ----------------
$ws.Range("B115").Select()
